# Updates cryptos list pricing/volume data (and fixes the Frax/TrustWalletToken
# row ordering) to match the latest scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.933.68'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '1.810.96'
$ws.Range("E3").Value = '  +1.49%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '''310.18'
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("D7").Value = '''0.4987'
$ws.Range("E7").Value = '  -2.61%  '
$ws.Range("D8").Value = '''0.3915'
$ws.Range("E8").Value = '  +1.60%  '
$ws.Range("D9").Value = '''0.09802'
$ws.Range("E9").Value = '  +25.38%  '
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("D11").Value = '''40.81'
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = '''6.413'
$ws.Range("E12").Value = '  +3.40%  '
$ws.Range("D13").Value = '''1.000'
$ws.Range("E13").Value = '  -0.13%  '
$ws.Range("D14").Value = '''20.40'
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("D15").Value = '1.810.43'
$ws.Range("E15").Value = '  +1.63%  '
$ws.Range("D16").Value = '''7.264'
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").Value = '''0.00001136'
$ws.Range("E17").Value = '  +5.59%  '
$ws.Range("D18").Value = '''92.22'
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("D19").Value = '''0.06635'
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").Value = '''17.15'
$ws.Range("E21").Value = '  +0.74%  '
$ws.Range("D22").Value = '''5.902'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").Value = '27.990.80'
$ws.Range("E23").Value = '  +0.55%  '
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("E25").Value = '  +0.53%  '
$ws.Range("D26").Value = '''158.61'
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D27").Value = '2.019.87'
$ws.Range("E27").Value = '  +1.63%  '
$ws.Range("D28").Value = '''20.50'
$ws.Range("E28").Value = '  +1.45%  '
$ws.Range("D29").Value = '''2.384'
$ws.Range("E29").Value = '  +0.74%  '
$ws.Range("D30").Value = '''126.68'
$ws.Range("E30").Value = '  +2.55%  '
$ws.Range("D31").Value = '''0.1064'
$ws.Range("E31").Value = '  -0.92%  '
$ws.Range("D32").Value = '''1.031'
$ws.Range("E32").Value = '  -0.45%  '
$ws.Range("D33").Value = '''5.544'
$ws.Range("E33").Value = '  +1.08%  '
$ws.Range("D34").Value = '''3.604'
$ws.Range("E34").Value = '  -0.82%  '
$ws.Range("D35").Value = '''0.06729'
$ws.Range("E35").Value = '  -5.06%  '
$ws.Range("D37").Value = '''8.840'
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("D38").Value = '''0.2134'
$ws.Range("E38").Value = '  +0.56%  '
$ws.Range("E39").Value = '  -1.57%  '
$ws.Range("E40").Value = '  -2.01%  '
$ws.Range("D41").Value = '''0.6158'
$ws.Range("E41").Value = '  +1.13%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''1.172'
$ws.Range("E42").Value = '  +1.55%  '
$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").Value = '''1.0000'
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").Value = '''13.10'
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("D46").Value = '''1.286'
$ws.Range("E46").Value = '  -2.46%  '
$ws.Range("D47").Value = '''3.689'
$ws.Range("E47").Value = '  -0.40%  '
$ws.Range("D48").Value = '''123.44'
$ws.Range("E48").Value = '  -2.02%  '
$ws.Range("D49").Value = '''1.933'
$ws.Range("E49").Value = '  +1.91%  '
$ws.Range("D50").Value = '''1.175'
$ws.Range("E50").Value = '  -2.12%  '
$ws.Range("D51").Value = '''0.06765'
$ws.Range("E51").Value = '  -1.20%  '
